$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the whole target range to be treated as Text before assigning values,
# so numeric-looking strings (e.g. "99", "157.14") are stored as text, matching source data.
$ws.Range("A1:K8").NumberFormat = "@"

$ws.Range("A1").Value = "venue"
$ws.Range("B1").Value = "date"
$ws.Range("C1").Value = "result"
$ws.Range("D1").Value = "ownTeam"
$ws.Range("E1").Value = "oppTeam"
$ws.Range("F1").Value = "batsman"
$ws.Range("G1").Value = "totalRuns"
$ws.Range("H1").Value = "totalBalls"
$ws.Range("I1").Value = "total4s"
$ws.Range("J1").Value = "total6s"
$ws.Range("K1").Value = "sr"

$ws.Range("A2").Value = " Abu Dhabi"
$ws.Range("B2").Value = " October 30 2020"
$ws.Range("C2").Value = "Royals won by 7 wickets (with 15 balls remaining)"
$ws.Range("D2").Value = "Kings XI Punjab"
$ws.Range("E2").Value = "Rajasthan Royals"
$ws.Range("F2").Value = "Chris Gayle  "
$ws.Range("G2").Value = "99"
$ws.Range("H2").Value = "63"
$ws.Range("I2").Value = "6"
$ws.Range("J2").Value = "8"
$ws.Range("K2").Value = "157.14"

$ws.Range("A3").Value = " Dubai (DSC)"
$ws.Range("B3").Value = " October 24 2020"
$ws.Range("C3").Value = "Kings XI won by 12 runs"
$ws.Range("D3").Value = "Kings XI Punjab"
$ws.Range("E3").Value = "Sunrisers Hyderabad"
$ws.Range("F3").Value = "Chris Gayle  "
$ws.Range("G3").Value = "20"
$ws.Range("H3").Value = "20"
$ws.Range("I3").Value = "2"
$ws.Range("J3").Value = "1"
$ws.Range("K3").Value = "100.00"

$ws.Range("A4").Value = " Sharjah"
$ws.Range("B4").Value = " October 15 2020"
$ws.Range("C4").Value = "Kings XI won by 8 wickets"
$ws.Range("D4").Value = "Kings XI Punjab"
$ws.Range("E4").Value = "Royal Challengers Bangalore"
$ws.Range("F4").Value = "Chris Gayle  "
$ws.Range("G4").Value = "53"
$ws.Range("H4").Value = "45"
$ws.Range("I4").Value = "1"
$ws.Range("J4").Value = "5"
$ws.Range("K4").Value = "117.77"

$ws.Range("A5").Value = " Abu Dhabi"
$ws.Range("B5").Value = " November 01 2020"
$ws.Range("C5").Value = "Super Kings won by 9 wickets (with 7 balls remaining)"
$ws.Range("D5").Value = "Kings XI Punjab"
$ws.Range("E5").Value = "Chennai Super Kings"
$ws.Range("F5").Value = "Chris Gayle  "
$ws.Range("G5").Value = "12"
$ws.Range("H5").Value = "19"
$ws.Range("I5").Value = "0"
$ws.Range("J5").Value = "0"
$ws.Range("K5").Value = "63.15"

$ws.Range("A6").Value = " Dubai (DSC)"
$ws.Range("B6").Value = " October 20 2020"
$ws.Range("C6").Value = "Kings XI won by 5 wickets (with 6 balls remaining)"
$ws.Range("D6").Value = "Kings XI Punjab"
$ws.Range("E6").Value = "Delhi Capitals"
$ws.Range("F6").Value = "Chris Gayle  "
$ws.Range("G6").Value = "29"
$ws.Range("H6").Value = "13"
$ws.Range("I6").Value = "3"
$ws.Range("J6").Value = "2"
$ws.Range("K6").Value = "223.07"

$ws.Range("A7").Value = " Sharjah"
$ws.Range("B7").Value = " October 26 2020"
$ws.Range("C7").Value = "Kings XI won by 8 wickets (with 7 balls remaining)"
$ws.Range("D7").Value = "Kings XI Punjab"
$ws.Range("E7").Value = "Kolkata Knight Riders"
$ws.Range("F7").Value = "Chris Gayle  "
$ws.Range("G7").Value = "51"
$ws.Range("H7").Value = "29"
$ws.Range("I7").Value = "2"
$ws.Range("J7").Value = "5"
$ws.Range("K7").Value = "175.86"

$ws.Range("A8").Value = " Dubai (DSC)"
$ws.Range("B8").Value = " October 18 2020"
$ws.Range("C8").Value = "Match tied (Kings XI won the one-over eliminator)"
$ws.Range("D8").Value = "Kings XI Punjab"
$ws.Range("E8").Value = "Mumbai Indians"
$ws.Range("F8").Value = "Chris Gayle  "
$ws.Range("G8").Value = "24"
$ws.Range("H8").Value = "21"
$ws.Range("I8").Value = "1"
$ws.Range("J8").Value = "2"
$ws.Range("K8").Value = "114.28"
